# Add 2022-Q4 data
#
# Strategy: the existing "2022-Q3" detail sheet becomes the new "2022-Q4"
# detail sheet (same underlying sheet object / sheetId), while a duplicate
# of it is created to preserve the old "2022-Q3" data under its original
# name (this matches the sheetId layout seen after the edit: 总计=1,
# 2022-Q4=2, 2022-Q3=3).

$wb = $excel.ActiveWorkbook

$wsQ3 = $wb.Worksheets.Item("2022-Q3")

# Duplicate the Q3 sheet right after itself; this duplicate will keep the
# old Q3 data and the name "2022-Q3".
$wsQ3.Copy($null, $wsQ3)

$wsQ4 = $wb.Worksheets.Item("2022-Q3")
$wsQ3Copy = $wb.Worksheets.Item("2022-Q3 (2)")

$wsQ4.Name = "2022-Q4"
$wsQ3Copy.Name = "2022-Q3"

# ---------------------------------------------------------------------
# Helper: write a value that must be stored as literal TEXT even when it
# looks like a number (Excel would otherwise auto-convert "1.03" -> 1.03).
# Uses a throwaway TEXT() formula then collapses it to a static value via
# copy / paste-values, which avoids creating any new number-format style.
# ---------------------------------------------------------------------
function Set-TextValue {
    param($ws, $addr, [string]$text)
    $escaped = $text.Replace('"', '""')
    $ws.Range($addr).Formula = '="' + $escaped + '"'
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

# ---------------------------------------------------------------------
# Rebuild the "2022-Q4" detail sheet with the new quarter's holdings.
# ---------------------------------------------------------------------

# Drop the old 5th data row (old sheet had 4 funds, new one only has 3).
$wsQ4.Rows.Item(5).Delete()

# Row 2 - fund 860006
Set-TextValue $wsQ4 "B2" "860006"
Set-TextValue $wsQ4 "C2" "光大阳光优选一年持有混合A"
Set-TextValue $wsQ4 "D2" "1.03"
Set-TextValue $wsQ4 "E2" "80.15"
Set-TextValue $wsQ4 "F2" "1.53"
Set-TextValue $wsQ4 "G2" "0.0158"
$wsQ4.Range("H2").Value = 9

# Row 3 - fund 860055
Set-TextValue $wsQ4 "B3" "860055"
Set-TextValue $wsQ4 "C3" "光大阳光优选一年持有混合B"
Set-TextValue $wsQ4 "D3" "0.00"
Set-TextValue $wsQ4 "E3" "80.15"
Set-TextValue $wsQ4 "F3" "1.53"
$wsQ4.Range("G3").Value = 0
$wsQ4.Range("H3").Value = 9

# Row 4 - fund 860056
Set-TextValue $wsQ4 "B4" "860056"
Set-TextValue $wsQ4 "C4" "光大阳光优选一年持有混合C"
Set-TextValue $wsQ4 "D4" "0.00"
Set-TextValue $wsQ4 "E4" "80.15"
Set-TextValue $wsQ4 "F4" "1.53"
$wsQ4.Range("G4").Value = 0
$wsQ4.Range("H4").Value = 9

# Match the "总计" sheet's header / first-column styling (style index 2)
# on the rebuilt sheet, replacing the old detail-sheet styling (index 1).
$wsTotal = $wb.Worksheets.Item("总计")
$wsTotal.Range("B1").Copy()
$wsQ4.Range("B1:H1").PasteSpecial(-4122)
$wsTotal.Range("A2").Copy()
$wsQ4.Range("A2:A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Update the "总计" (summary) sheet: existing row 2 becomes the new
# 2022-Q4 summary, and a new row 3 is appended with the old 2022-Q3
# summary data that used to live in row 2.
# ---------------------------------------------------------------------

$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)
$wsTotal.Range("A3").Value = 1
Set-TextValue $wsTotal "B3" "2022-Q3"
$wsTotal.Range("C3").Value = 4
$wsTotal.Range("D3").Value = 0.44

Set-TextValue $wsTotal "B2" "2022-Q4"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.02

$excel.CutCopyMode = $false
